$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-02-03 19:01:12"

for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
